$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Intensity Distribution")

$ws.Range("B12").Value = 0.061110198497772217
$ws.Range("C12").Value = 0.12222039699554443
$ws.Range("B13").Value = 0.39985013008117676
$ws.Range("C13").Value = 0.76104438304901123
$ws.Range("B14").Value = 1.0496139526367188
$ws.Range("C14").Value = 1.5414133071899414
$ws.Range("B15").Value = 1.9699152708053589
$ws.Range("C15").Value = 1.9327237606048584
$ws.Range("B16").Value = 2.8891737461090088
$ws.Range("C16").Value = 2.0573775768280029
$ws.Range("B17").Value = 3.5554912090301514
$ws.Range("C17").Value = 2.7312109470367432
$ws.Range("B18").Value = 3.8809151649475098
$ws.Range("C18").Value = 3.875849723815918
$ws.Range("B19").Value = 3.9284200668334961
$ws.Range("C19").Value = 4.6171903610229492
$ws.Range("B20").Value = 3.7340297698974609
$ws.Range("C20").Value = 4.5281534194946289
$ws.Range("B21").Value = 3.0919084548950195
$ws.Range("C21").Value = 3.9345018863677979
$ws.Range("B22").Value = 2.1876351833343506
$ws.Range("C22").Value = 2.99216628074646
$ws.Range("B23").Value = 1.2634670734405518
$ws.Range("C23").Value = 1.9267702102661133
$ws.Range("B24").Value = 0.54180580377578735
$ws.Range("C24").Value = 0.96550935506820679
$ws.Range("B25").Value = 0.14181087911128998
$ws.Range("C25").Value = 0.28362175822257996
$ws.Range("B32").Value = 0.0056865592487156391
$ws.Range("C32").Value = 0.011373118497431278
$ws.Range("B33").Value = 0.055282726883888245
$ws.Range("C33").Value = 0.11056545376777649
$ws.Range("B34").Value = 0.14471666514873505
$ws.Range("C34").Value = 0.26228576898574829
$ws.Range("B35").Value = 0.24896632134914398
$ws.Range("C35").Value = 0.38894027471542358
$ws.Range("B36").Value = 0.32041388750076294
$ws.Range("C36").Value = 0.43642884492874146
$ws.Range("B37").Value = 0.32653334736824036
$ws.Range("C37").Value = 0.39473026990890503
$ws.Range("B38").Value = 0.26963251829147339
$ws.Range("C38").Value = 0.31353896856307983
$ws.Range("B39").Value = 0.18607331812381744
$ws.Range("C39").Value = 0.26941037178039551
$ws.Range("B40").Value = 0.12421361356973648
$ws.Range("C40").Value = 0.24842722713947296
$ws.Range("B41").Value = 0.10216101258993149
$ws.Range("C41").Value = 0.20432202517986298
$ws.Range("B42").Value = 0.088003858923912048
$ws.Range("C42").Value = 0.1760077178478241
$ws.Range("B43").Value = 0.26527580618858337
$ws.Range("C43").Value = 0.33911675214767456
$ws.Range("B44").Value = 1.1317200660705566
$ws.Range("C44").Value = 2.0022642612457275
$ws.Range("B45").Value = 2.3193156719207764
$ws.Range("C45").Value = 4.3246865272521973
$ws.Range("B46").Value = 3.20468807220459
$ws.Range("C46").Value = 6.0618529319763184
$ws.Range("B47").Value = 3.3346977233886719
$ws.Range("C47").Value = 6.3328094482421875
$ws.Range("B48").Value = 2.6490767002105713
$ws.Range("C48").Value = 5.0270686149597168
$ws.Range("B49").Value = 1.496807336807251
$ws.Range("C49").Value = 2.827465295791626
$ws.Range("B50").Value = 0.76980262994766235
$ws.Range("C50").Value = 0.86034637689590454
$ws.Range("B51").Value = 1.6585675477981567
$ws.Range("C51").Value = 3.3171350955963135
$ws.Range("B52").Value = 3.6877670288085938
$ws.Range("C52").Value = 7.3755340576171875
$ws.Range("B53").Value = 5.5730428695678711
$ws.Range("C53").Value = 11.13145637512207
$ws.Range("B54").Value = 6.6729903221130371
$ws.Range("C54").Value = 12.895022392272949
$ws.Range("B55").Value = 6.5645241737365723
$ws.Range("C55").Value = 11.912178993225098
$ws.Range("B56").Value = 5.2994780540466309
$ws.Range("C56").Value = 8.63383960723877
$ws.Range("B57").Value = 3.4213066101074219
$ws.Range("C57").Value = 4.5907869338989258
$ws.Range("B58").Value = 1.7351020574569702
$ws.Range("C58").Value = 2.0638277530670166
$ws.Range("B59").Value = 0.85307943820953369
$ws.Range("C59").Value = 1.7061588764190674
$ws.Range("B60").Value = 0.56355828046798706
$ws.Range("C60").Value = 1.1271165609359741
$ws.Range("B61").Value = 0.26289764046669006
$ws.Range("C61").Value = 0.52579528093338013
$ws.Range("B62").Value = 0.060074236243963242
$ws.Range("C62").Value = 0.12014847248792648
$ws.Range("B63").Value = 0.52874898910522461
$ws.Range("C63").Value = 1.0574979782104492
$ws.Range("B64").Value = 1.8271431922912598
$ws.Range("C64").Value = 3.6542863845825195
$ws.Range("B65").Value = 3.4232687950134277
$ws.Range("C65").Value = 6.8465375900268555
$ws.Range("B66").Value = 4.6186394691467285
$ws.Range("C66").Value = 9.237278938293457
$ws.Range("B67").Value = 4.9032154083251953
$ws.Range("C67").Value = 9.7286376953125
$ws.Range("B68").Value = 4.2315406799316406
$ws.Range("C68").Value = 7.99702787399292
$ws.Range("B69").Value = 2.987107515335083
$ws.Range("C69").Value = 4.7305364608764648
$ws.Range("B70").Value = 1.8217302560806274
$ws.Range("C70").Value = 1.6008585691452026
$ws.Range("B71").Value = 2.1493856906890869
$ws.Range("C71").Value = 1.5313806533813477
$ws.Range("B72").Value = 3.8377039432525635
$ws.Range("C72").Value = 3.9504415988922119
